$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.443.14'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.900.83'
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''602.49'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '''167.80'
$ws.Range("E6").Value = '  +2.46%  '
$ws.Range("D7").Value = '3.899.87'
$ws.Range("E7").Value = '  +2.50%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").Value = '''0.168'
$ws.Range("E10").Value = '  -1.46%  '
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").Value = '''0.461'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '''37.40'
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("D15").Value = '4.550.35'
$ws.Range("E15").Value = '  +2.39%  '
$ws.Range("D16").Value = '3.920.06'
$ws.Range("E16").Value = '  +3.07%  '
$ws.Range("D17").Value = '68.502.86'
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = '''7.48'
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("D19").Value = '''17.33'
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").Value = '''11.10'
$ws.Range("E21").Value = '  -3.22%  '
$ws.Range("D22").Value = '''490.47'
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D23").Value = '''0.727'
$ws.Range("E23").Value = '  +0.86%  '
$ws.Range("E24").Value = '  +3.97%  '
$ws.Range("D25").Value = '''84.70'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").Value = '''2.24'
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("D31").Value = '4.051.02'
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("D32").Value = '''2.38'
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").Value = '''7.74'
$ws.Range("E33").Value = '  -3.40%  '
$ws.Range("D34").Value = '''31.79'
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '3.853.32'
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("D39").Value = '''5.93'
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("D40").Value = '''3.18'
$ws.Range("E40").Value = '  +4.33%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("D43").Value = '''431.56'
$ws.Range("E43").Value = '  +0.47%  '
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '''48.07'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").Value = '''8.55'
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '''0.000272'
$ws.Range("E48").Value = '  +17.89%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''142.75'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").Value = '2.801.50'
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("D51").Value = '''39.49'
$ws.Range("E51").Value = '  -0.19%  '
